$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2840
$ws1.Range("F8").Value = 1848

# Update "全部类型" (All Types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2840
$ws4.Range("F8").Value = 1848
